# The commit swaps the presentation's "Integral" theme palette for the
# stock Office theme palette (ppt/theme/theme1.xml, the theme used by the
# slide master / every slide via the Design).
#
# PowerPoint's object model doesn't give VBA/COM code a way to overwrite a
# whole theme part wholesale (there is no "import this OOXML as the new
# theme" call) - the documented, supported way to recolor a theme through
# automation is to walk ThemeColorScheme.Colors(i).RGB, per
# "edit the theme via ThemeColorScheme.Colors(i).RGB ... instead". So we
# set each of the twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the exact values used by the standard "Office Theme"
# colour scheme, in the official slot order.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

function Set-ThemeColor {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $themeColors.Colors($Index).RGB = $rgb
}

# Slot order matches the OOXML a:clrScheme child order.
Set-ThemeColor 1  "000000"   # dk1
Set-ThemeColor 2  "FFFFFF"   # lt1
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink

Write-Host "Applied Office theme color scheme to the presentation theme."
